$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

Set-TextValue "D2" "26.280.67"
Set-TextValue "E2" "  -0.49%  "
Set-TextValue "D3" "1.667.47"
Set-TextValue "E3" "  -0.55%  "
Set-TextValue "D4" "1.005"
Set-TextValue "E4" "  +0.35%  "
Set-TextValue "D5" "219.25"
Set-TextValue "E5" "  +1.24%  "
Set-TextValue "D6" "0.5242"
Set-TextValue "E6" "  -1.18%  "
Set-TextValue "D7" "1.005"
Set-TextValue "E7" "  +0.33%  "
Set-TextValue "D8" "0.2711"
Set-TextValue "E8" "  +0.29%  "
Set-TextValue "D9" "0.06337"
Set-TextValue "E9" "  -1.00%  "
Set-TextValue "D10" "21.01"
Set-TextValue "E10" "  -3.14%  "
Set-TextValue "D11" "0.07757"
Set-TextValue "E11" "  -0.62%  "
Set-TextValue "B12" "Polkadot"
Set-TextValue "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "4.452"
Set-TextValue "E12" "  -1.20%  "
Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.657.69"
Set-TextValue "E13" "  -1.08%  "
Set-TextValue "D14" "1.893.25"
Set-TextValue "E14" "  -0.65%  "
Set-TextValue "D15" "0.5473"
Set-TextValue "E15" "  -1.49%  "
Set-TextValue "D16" "0.0₅8238"
Set-TextValue "E16" "  -0.88%  "
Set-TextValue "D17" "65.00"
Set-TextValue "E17" "  -0.83%  "
Set-TextValue "D18" "26.294.46"
Set-TextValue "E18" "  -0.63%  "
Set-TextValue "E19" "  +0.39%  "
Set-TextValue "D20" "4.669"
Set-TextValue "E20" "  -1.37%  "
Set-TextValue "D21" "195.28"
Set-TextValue "E21" "  +0.86%  "
Set-TextValue "D22" "10.16"
Set-TextValue "E22" "  -1.47%  "
Set-TextValue "D23" "6.088"
Set-TextValue "E23" "  -3.96%  "
Set-TextValue "D24" "1.007"
Set-TextValue "E24" "  +0.46%  "
Set-TextValue "D25" "139.90"
Set-TextValue "E25" "  -1.81%  "
Set-TextValue "D26" "0.1241"
Set-TextValue "E26" "  -3.28%  "
Set-TextValue "D27" "7.206"
Set-TextValue "E27" "  -2.72%  "
Set-TextValue "D28" "16.18"
Set-TextValue "E28" "  -0.47%  "
Set-TextValue "E29" "  -0.51%  "
Set-TextValue "D30" "0.06154"
Set-TextValue "E30" "  -1.51%  "
Set-TextValue "E31" "  +1.00%  "
Set-TextValue "D32" "3.594"
Set-TextValue "E32" "  -0.39%  "
Set-TextValue "D33" "3.294"
Set-TextValue "E33" "  -4.38%  "
Set-TextValue "D34" "1.637"
Set-TextValue "E34" "  -2.35%  "
Set-TextValue "D35" "0.9761"
Set-TextValue "E35" "  -3.01%  "
Set-TextValue "B36" "MXToken"
Set-TextValue "C36" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D36" "2.791"
Set-TextValue "E36" "  +0.42%  "
Set-TextValue "B37" "HuobiToken"
Set-TextValue "C37" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D37" "2.419"
Set-TextValue "E37" "  -0.38%  "
Set-TextValue "D38" "0.5731"
Set-TextValue "E38" "  -5.39%  "
Set-TextValue "D39" "0.01608"
Set-TextValue "E39" "  -1.43%  "
Set-TextValue "D40" "6.042"
Set-TextValue "E40" "  -2.41%  "
Set-TextValue "D41" "0.8569"
Set-TextValue "E41" "  -0.84%  "
Set-TextValue "E42" "  +0.41%  "
Set-TextValue "D43" "1.023.92"
Set-TextValue "E43" "  -5.62%  "
Set-TextValue "D44" "100.44"
Set-TextValue "E44" "  +0.30%  "
Set-TextValue "D45" "1.809.15"
Set-TextValue "E45" "  -0.68%  "
Set-TextValue "D46" "58.06"
Set-TextValue "E46" "  +1.77%  "
Set-TextValue "D47" "0.0₈109"
Set-TextValue "E47" "  -0.56%  "
Set-TextValue "D48" "1.009"
Set-TextValue "E48" "  +0.82%  "
Set-TextValue "D49" "8.080"
Set-TextValue "E49" "  -0.54%  "
Set-TextValue "D50" "1.489"
Set-TextValue "E50" "  +1.72%  "
Set-TextValue "D51" "0.05186"
Set-TextValue "E51" "  -0.36%  "
